$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (2..123) down to (3..124) by copying each row's
# values/formatting into the row below it, working from the bottom up so that
# no data is overwritten before it is copied. This avoids the format bleed
# that Rows.Insert() causes when inserting directly below the bold header row.
for ($r = 123; $r -ge 2; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $src.Copy($dst)
}

# Update the sheet dimension to reflect the extra row.
$ws.Range("A1:R124").Select() | Out-Null

# Populate row 2 with the new weekly record.
$ws.Cells.Item(2, 1).Value = 10
$ws.Cells.Item(2, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(2, 3).Value = "La Araucanía"
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 4).Value = 44812
$ws.Cells.Item(2, 5).Value = 9
$ws.Cells.Item(2, 6).Value = 100114002
$ws.Cells.Item(2, 7).Value = "Camote"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 50
$ws.Cells.Item(2, 11).Value = 20000
$ws.Cells.Item(2, 12).Value = 20000
$ws.Cells.Item(2, 13).Value = 20000
$ws.Cells.Item(2, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(2, 15).Value = "Perú"
$ws.Cells.Item(2, 16).Value = 1000
$ws.Cells.Item(2, 17).Value = 20
$ws.Cells.Item(2, 18).Value = "Hortaliza"
